# ToDoItems.xlsx: rework the sample fixture data sheet so that the
# ExcelFixtureRowHandler "previousRow" feature has something to show off:
#   - columns are reordered from (description, category, subcategory) to
#     (category, subcategory, description)
#   - repeated category/subcategory values (relative to the row above) are
#     suppressed: category is blanked (but keeps the "continuation" style),
#     subcategory is removed outright when unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$lastRow = 13

# --- 1. Rotate the A/B/C columns: new A = old B (category), new B = old C
#        (subcategory), new C = old A (description). ----------------------
for ($r = 1; $r -le $lastRow; $r++) {
    $oldA = $ws.Cells.Item($r, 1).Value2
    $oldB = $ws.Cells.Item($r, 2).Value2
    $oldC = $ws.Cells.Item($r, 3).Value2

    $ws.Cells.Item($r, 1).Value2 = $oldB
    $ws.Cells.Item($r, 2).Value2 = $oldC
    $ws.Cells.Item($r, 3).Value2 = $oldA
}

# --- 2. Give every data row's category cell (col A) the "continuation"
#        style (same style already used for the old repeated category
#        cells) by copying it from a cell that already carries it. --------
$ws.Cells.Item(3, 2).Copy() | Out-Null
for ($r = 3; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- 3. The subcategory column (col B) no longer needs the continuation
#        style at all now that repeats are removed outright - reset it
#        back to the Normal style first. -----------------------------------
for ($r = 3; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Style = "Normal"
}

# --- 4. Blank out (but keep styled) the category cell wherever it is the
#        same as the row above; remove the subcategory cell entirely
#        wherever it is the same as the row above. -------------------------
$prevCategory = $ws.Cells.Item(2, 1).Value2
$prevSubcategory = $ws.Cells.Item(2, 2).Value2

for ($r = 3; $r -le $lastRow; $r++) {
    $thisCategory = $ws.Cells.Item($r, 1).Value2
    $thisSubcategory = $ws.Cells.Item($r, 2).Value2

    if ($thisCategory -eq $prevCategory) {
        $ws.Cells.Item($r, 1).Value2 = $null
    } else {
        $prevCategory = $thisCategory
    }

    if ($thisSubcategory -eq $prevSubcategory) {
        $ws.Cells.Item($r, 2).Clear() | Out-Null
    } else {
        $prevSubcategory = $thisSubcategory
    }
}

# --- 5. Swap the column widths around to follow the data (A<-B, B<-C, C<-A)
$ws.Columns.Item(1).ColumnWidth = 11.6
$ws.Columns.Item(2).ColumnWidth = 10.8
$ws.Columns.Item(3).ColumnWidth = 24.3

# --- 6. Freeze the first two columns (instead of just the first) and
#        refresh the sheet's remembered selections to match. ---------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D2:D4").Select() | Out-Null
